# Add the "Julia" function-translation column (F) to the export sheet,
# and fix a stale lookup in C20 (subtract(x,y) should show "x-y", not "x/y").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- fix pre-existing typo: C20 pointed at the wrong shared string ---
$ws.Range("C20").Value2 = "x-y"

# --- new column F header + values, formatted like the matching E-column cell ---
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value2 = "Julia"
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value2 = "abs(x)"
$ws.Range("E4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value2 = " +(x,y)"
$ws.Range("E3").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value2 = "ceil(x)"
$ws.Range("E4").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value2 = "^(x,3)"
$ws.Range("E4").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value2 = " /(x,y)"
$ws.Range("E3").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value2 = "exp(x)"
$ws.Range("E3").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value2 = "floor(x)"
$ws.Range("E4").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value2 = "log(x)"
$ws.Range("E3").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value2 = "log(x)"
$ws.Range("E4").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value2 = "log(base,x)"
$ws.Range("E3").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value2 = "log10(x)"
$ws.Range("E3").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value2 = "log2(x)"
$ws.Range("E4").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value2 = "*(x,y)"
$ws.Range("E4").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16").Value2 = "^(x,y)"
$ws.Range("E3").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value2 = "sign(x)"
$ws.Range("E3").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value2 = "sqrt(x)"
$ws.Range("E4").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value2 = "^(x,2)"
$ws.Range("E4").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value2 = " -(x,y)"
$ws.Range("E3").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F21").Value2 = "max(x, y)"
$ws.Range("E3").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value2 = "max(x, y, z)"
$ws.Range("E3").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("F23").Value2 = "min(x, y)"
$ws.Range("E3").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("F24").Value2 = "min(x, y, z)"
$ws.Range("E3").Copy()
$ws.Range("F25").PasteSpecial(-4122)
$ws.Range("F25").Value2 = "x+y"
$ws.Range("E3").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value2 = "x-y"
$ws.Range("E3").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value2 = "x*y"
$ws.Range("E3").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value2 = "x/y"
$ws.Range("E3").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value2 = "x^y"
$ws.Range("E4").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value2 = "x-y > 0 ? 1 : 2"
$ws.Range("E4").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value2 = "x-y == 0 ? 1 : 2"
$ws.Range("E4").Copy()
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("F34").Value2 = "x-y >= 0 ? 1 : 2"

$excel.CutCopyMode = $false

# --- column widths: D widened, new column F gets its own width ---
$ws.Columns.Item(4).ColumnWidth = 24.3
$ws.Columns.Item(6).ColumnWidth = 13.0

# --- view: scroll down and select H40, matching the edited workbook ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("H40").Select()
